# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (strings that cannot be mistaken for numbers)
$ws.Range("D2").Value = '65.426.07'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = '3.285.46'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.281.42'
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("E10").Value = '  -6.35%  '
$ws.Range("E11").Value = '  -3.58%  '
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '3.802.59'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("E16").Value = '  -3.82%  '
$ws.Range("D17").Value = '65.520.80'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '3.283.89'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("E24").Value = '  -3.41%  '
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("E26").Value = '  -0.38%  '
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("E28").Value = '  -3.62%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("E30").Value = '  -4.32%  '
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E33").Value = '  -10.74%  '
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("D35").Value = '3.807.36'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  -3.55%  '
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("E40").Value = '  +6.26%  '
$ws.Range("E41").Value = '  -5.89%  '
$ws.Range("E42").Value = '  -6.02%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0673'
$ws.Range("E43").Value = '  -8.24%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E44").Value = '  -6.86%  '
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("E46").Value = '  -4.38%  '
$ws.Range("E47").Value = '  -5.31%  '
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("E50").Value = '  -4.08%  '
$ws.Range("E51").Value = '  +5.44%  '

# Numeric-looking Price strings must be forced to remain text,
# matching the original inline-string formatted values (e.g. "0.998").
# Temporarily mark the cell as Text, assign the value, then clear the
# formatting override so no stray style index is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.34'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.84'
$ws.Range("D6").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.568'
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.567'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.28'
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '628.86'
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.34'
$ws.Range("D16").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.58'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.86'
$ws.Range("D21").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.96'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.58'
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.95'
$ws.Range("D25").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.71'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.32'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.54'
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.34'
$ws.Range("D30").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '557.93'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.64'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.78'
$ws.Range("D34").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.77'
$ws.Range("D38").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.24'
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.12'
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("D44").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.05'
$ws.Range("D47").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.126'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.49'
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.24'
$ws.Range("D51").ClearFormats()
